$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "navegate" -> "navigate"
$ws.Range("B2").Value = "navigate"

# Move active cell selection to C2
$ws.Range("C2").Select()
